# This workbook tracks weekly "Piña" (pineapple) price observations for the
# "Vega Modelo de Temuco" market. Each new weekly update inserts a fresh row
# at the top of the data block (row 706) carrying the same catalog/lookup
# values as the most recent prior observation (row 706) but stamped with a
# new date. All subsequent rows shift down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 706, pushing the existing data (rows 706-808)
# down to rows 707-809.
$ws.Rows.Item(706).Insert()

# The row that used to be 706 is now at 707. Duplicate its full contents
# (values, number formats, styles) into the newly inserted row 706.
$ws.Range("A707:T707").Copy($ws.Range("A706:T706"))

# Stamp the new row with this week's date (serial 45180 = 2023-09-11),
# while keeping every other column (grade, volumes, prices, origin, etc.)
# identical to the prior observation that was just copied down.
$ws.Cells.Item(706, 4).Value = 45180
